$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NAME column (A) with the new names
$ws.Range("A4").Value = "Edutech  Christain"
$ws.Range("A3").Value = "Onyekachi  Ekenechukwu"
$ws.Range("A2").Value = "Chinedu  Patrick"

# Update the selected cell in the sheet view
$ws.Range("A8").Select()
